$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# The rows describing "REPOSICAO DE BLOQUETE" / "REPOSICAO DE BLOQUETE INV"
# that used to live on their own sheet ("n10") are being folded into the
# "reposicao" sheet instead, and the now-empty "n10" sheet is removed.

$reposicao = $wb.Worksheets.Item("reposicao")
$reposicao.Activate()

$reposicao.Cells.Item(16, 1).Value = "738000"
$reposicao.Cells.Item(16, 2).Value = "REPOSIÇÃO DE BLOQUETE"
$reposicao.Cells.Item(16, 3).Value = "Reposicao"

$reposicao.Cells.Item(17, 1).Value = "740000"
$reposicao.Cells.Item(17, 2).Value = "REPOSIÇÃO DE BLOQUETE INV"
$reposicao.Cells.Item(17, 3).Value = "Reposicao"

# Match the formatting the two rows carried on the original "n10" sheet
# (left-aligned text style used throughout column A of these lookup sheets).
$colA = $reposicao.Range("A16:A17")
$colA.NumberFormat = "@"
$colA.HorizontalAlignment = -4131
$colA.Font.Name = "Calibri"
$colA.Font.Size = 11

$reposicao.Range("B24").Select()

# Remove the sheet whose rows were just merged into "reposicao".
$n10 = $wb.Worksheets.Item("n10")
$n10.Delete()

# "n3" becomes the active tab after the reshuffle.
$n3 = $wb.Worksheets.Item("n3")
$n3.Activate()
